# Insert a new data row at row 40 (shifts existing rows 40-83 down to 41-84)
# and populate it with a new weekly "Orégano" price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("40:40").Insert()

$ws.Range("A40").Value = 9
$ws.Range("B40").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C40").Value = "Metropolitana"
$ws.Range("D40").Value = 44874
$ws.Range("E40").Value = 13
$ws.Range("F40").Value = 100112029
$ws.Range("G40").Value = "Orégano"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 8
$ws.Range("K40").Value = 15000
$ws.Range("L40").Value = 18000
$ws.Range("M40").Value = 16500
$ws.Range("N40").Value = "`$/docena de atados"
$ws.Range("O40").Value = "Región Metropolitana"
$ws.Range("P40").Value = 5500
$ws.Range("Q40").Value = 3
$ws.Range("R40").Value = "Hortaliza"
